$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# ---------------------------------------------------------------------
# Shape "Rectangle 4" (SR-MPLS probe-query figure box)
# ---------------------------------------------------------------------
$srMpls = $s.Shapes.Item(3)

# Move the box up slightly (y: 102581 -> 57150 EMU)
$srMpls.Top = 4.5

# Update the figure caption text (keep it a single run)
$srMplsText = $srMpls.TextFrame.TextRange
$srMplsCaption = $srMplsText.Paragraphs($srMplsText.Paragraphs().Count).Runs(1)
$srMplsCaption.Text = "     Figure: Example Probe Query Message for SR-MPLS Policy"

# ---------------------------------------------------------------------
# Shape "Rectangle 8" (SRv6 probe-query figure box)
# ---------------------------------------------------------------------
$srv6 = $s.Shapes.Item(6)
$srv6Text = $srv6.TextFrame.TextRange

# Insert a new "Next Header = 43 (Routing Header)" line right before the
# blank "." line that follows the Destination IP Address line.
$blankAfterDest = $srv6Text.Paragraphs(7).Runs(1)
$blankAfterDest.InsertBefore(".  Next Header = 43 (Routing Header)                            .`r") | Out-Null

# Update the SID List line text, then add a new "Next Header = 17 (UDP)"
# line right after it. (The paragraph insert above shifted this line from
# index 10 to index 11.)
$sidListRun = $srv6Text.Paragraphs(11).Runs(1)
$sidListRun.Text = ".  <SID List>                                                   ."
$sidListRun.InsertAfter("`r.  Next Header = 17 (UDP)                                       .") | Out-Null

# Update the figure caption text (keep it a single run)
$srv6Caption = $srv6Text.Paragraphs($srv6Text.Paragraphs().Count).Runs(1)
$srv6Caption.Text = "       Figure: Example Probe Query Message for SRv6 Policy"

# Reposition the box (the height grows automatically via the shape's
# auto-fit behaviour once the two new lines above are added).
$srv6.Top = 172.57968903937007
